# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, pushing the old N:Q ("Late"/heading/heading/Outstanding) block to
#   O:R. The new column inherits its width from the column to its left (M),
#   same as Excel does on a manual column insert.
# - Make "Repayment schedule" the active sheet/tab, with S7 selected
#   (previously "NewLoanInput" was the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column immediately to the left (M) before the
# insert shifts everything one column to the right.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column at N (existing N..Q shift right to O..R).
$ws.Columns("N").Insert()

# New column N picks up the width of the column to its left, like a real
# Excel "Insert" does.
$ws.Columns("N").ColumnWidth = $leftWidth

# Switch to this sheet and select S7, matching the saved view state.
$ws.Activate()
$ws.Range("S7").Select()
